$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 truth-table corrections ---------------------------------------
# K16: RegDst0 1 -> 0
$ws.Range("K16").Value = 0
# P16: MemtoReg0 1 -> 0
$ws.Range("P16").Value = 0
# Q16: RegDst1 0 -> 1
$ws.Range("Q16").Value = 1

# V16 ("Hex" column): drop the stray highlighted "0521" value, replace with
# the regular "0111" value/format used by the rest of the table (copy format
# from V2, which already carries the un-highlighted style).
$ws.Range("V2").Copy()
$ws.Range("V16").PasteSpecial(-4122)
$ws.Range("V16").Value = "0111"
$excel.CutCopyMode = 0

# --- Selection bookkeeping --------------------------------------------------
# Mirrors clicking the row-16 header while D16 was the active cell: the whole
# row becomes selected (sqref A16:XFD16) with D16 remembered as the anchor
# cell.
[void]$ws.Rows("16:16").Select()
try { $excel.ActiveCell = $ws.Range("D16") } catch {}
